$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 2) with the MCH148 collection record.
$ws.Range("A2").Value = "MCH148"
$ws.Range("C2").Value = "THE CONTRIBUTION OF THE EUROPEANS TO WORLD CIVILIZATION"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Format the new row with the same font used elsewhere on the sheet
# (Calibri 10pt, automatic/theme text color) and carry that formatting
# onto the otherwise-empty D2/H2 cells as well.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-establish the frozen header row / selection on the new data row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:I2").Select() | Out-Null
